$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210310339927673
$ws.Range("B1").Value = 2.516122341156006
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.739291548728943
$ws.Range("E1").Value = 1.162550449371338
